$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Update sheet 1 ("herzo_student") row 2 contents and remove row 3
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "test1a"
$ws1.Range("B2").Value = "test1b"
$ws1.Range("C2").Value = "test-ka"
$ws1.Range("F2").Formula = '=LOWER(D2&"."&E2&"@sbs-herzogenaurach.de")'
$ws1.Range("G2").Value = "1234abc!"

# Apply formatting (vertical-centered, black font) to A2:B2 the same way it was
# applied in the authored workbook: build the style on A2 first, then clone the
# resulting cell format onto B2 so only a single new style record is produced.
$a2 = $ws1.Range("A2")
$a2.VerticalAlignment = -4108
$a2.Font.Color = 0x000000
$a2.Copy() | Out-Null
$ws1.Range("B2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Remove the now-obsolete third data row
$ws1.Rows("3:3").Delete() | Out-Null

# Update the selection / page setup to match the authored file
$ws1.Range("G2").Select() | Out-Null
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 2) Add the new sheet "hoe_student" right after "herzo_student"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "hoe_student"

# Column widths matching the authored worksheet (closest achievable values)
$ws2.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 15.333333333333334
$ws2.Columns.Item(4).ColumnWidth = 14.833333333333334
$ws2.Columns.Item(6).ColumnWidth = 34.5
$ws2.Columns.Item(7).ColumnWidth = 23.333333333333332

# Header row (identical wording/styling to sheet 1's header row)
$ws2.Range("A1").Value = "givenname_raw"
$ws2.Range("B1").Value = "surname_raw"
$ws2.Range("C1").Value = "class"
$ws2.Range("D1").Value = "givenname"
$ws2.Range("E1").Value = "surname"
$ws2.Range("F1").Value = "UPN"
$ws2.Range("G1").Value = "password"

$ws2.Range("A1:C1").Font.Bold = $true
$ws2.Range("A1:C1").Interior.Color = 65535
$ws2.Range("D1:G1").Font.Bold = $true

# Data row
$ws2.Range("A2").Value = "kjhk test"
$ws2.Range("B2").Value = "lkjlkj"
$ws2.Range("D2").Formula = '=SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(A2,"Ä","Ae"),"Ö","Oe"),"Ü","Ue"),"ß","ss"),"ä","ae"),"ö","oe"),"ü","ue")," ","-")'
$ws2.Range("E2").Formula = '=SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(B2,"Ä","Ae"),"Ö","Oe"),"Ü","Ue"),"ß","ss"),"ä","ae"),"ö","oe"),"ü","ue")," ","-")'
$ws2.Range("F2").Formula = '=LOWER(D2&"."&E2&"@sbs-hoechstadt.de")'
$ws2.Range("G2").Value = "1234abc!"

# Apply the same vertical-centered black font styling to A2:B2
$a2b = $ws2.Range("A2")
$a2b.VerticalAlignment = -4108
$a2b.Font.Color = 0x000000
$a2b.Copy() | Out-Null
$ws2.Range("B2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws2.Range("G2").Select() | Out-Null
